{"js": "// Rename \"curve_id\" -> \"rqc_series_id\" and \"sample_amount\" -> \"relative_sample_amount\"\n// throughout the document body (both occurrences of each, in the two bullet\n// sentences that describe the new table columns).\n\nconst body = context.document.body;\n\nconst curveIdResults = body.search(\"curve_id\", { matchCase: true });\ncurveIdResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < curveIdResults.items.length; i++) {\n  curveIdResults.items[i].insertText(\"rqc_series_id\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconst sampleAmountResults = body.search(\"sample_amount\", { matchCase: true });\nsampleAmountResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < sampleAmountResults.items.length; i++) {\n  sampleAmountResults.items[i].insertText(\"relative_sample_amount\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Rename \"curve_id\" -> \"rqc_series_id\" and \"sample_amount\" -> \"relative_sample_amount\"\n# throughout the document body (both occurrences of each, in the two bullet\n# sentences that describe the new table columns).\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"curve_id\"\n$find.Replacement.Text = \"rqc_series_id\"\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"sample_amount\"\n$find2.Replacement.Text = \"relative_sample_amount\"\n$find2.Execute($find2.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n"}
